$d = $word.ActiveDocument

# 1) Insert a new empty "Normal" paragraph right before the
#    "Ver no Jupiter..." paragraph (i.e. right after the
#    "LOB1036: Geometria Analitica (Requisito fraco)" paragraph).
$find1 = $d.Content.Find
$find1.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false,
               $true, 1, $false, "", 0)
$find1.Parent.InsertParagraphBefore()

# 2) The "(c) 2020 . Contact: ..." paragraph becomes empty, and is
#    split into two paragraphs: the first stays plain "Normal", the
#    second carries PageBreakBefore + left alignment.
$find2 = $d.Content.Find
$find2.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
               $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyrightRange = $find2.Parent
$copyrightPara = $copyrightRange.Paragraphs(1)

# Split: add a new paragraph mark after this one, so we end up with
# two empty paragraphs where there used to be one with text.
$copyrightPara.Range.InsertParagraphAfter()

# Remove the text itself, leaving the (now first) paragraph empty.
$copyrightRange.Text = ""

# The second (new) paragraph gets the page break + left alignment
# that used to live on the following empty paragraph's sibling.
$secondPara = $copyrightPara.Next()
$secondPara.Format.PageBreakBefore = $true
$secondPara.Format.Alignment = 0

$d.Saved = $false
